$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dic_Disagg_Kategorien")

# Copy formatting from the row above (row 19) into the new row 20,
# then set the new values - mirrors adding a new data row to the table.
$ws.Range("A19:C19").Copy() | Out-Null
$ws.Range("A20:C20").PasteSpecial(-4122) | Out-Null

$ws.Range("A20").Value = "K_ZUORDN"
$ws.Range("B20").Value = "Zuordnung"
$ws.Range("C20").Value = "XXXZuordnung"
